$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new experiment record as row 32, based on the formatting of row 29
# (the most recent entry using the same model/dataset/preprocess settings).
$ws.Range("A29:Z29").Copy($ws.Range("A32:Z32"))

# Range.Copy() does not replicate cells that were blank in the source row,
# so re-copy those specifically to keep their (blank) styling.
$ws.Range("P29:Q29").Copy($ws.Range("P32:Q32"))
$ws.Range("Z29").Copy($ws.Range("Z32"))

# Row 29's "epoch" cell (L29) used a highlighted style; the new row should
# use the plain style instead, matching its neighboring cells.
$ws.Range("K29").Copy($ws.Range("L32"))

# Match row 29's height on the new row.
$ws.Rows(32).RowHeight = $ws.Rows(29).RowHeight

# Now fill in the values that differ for this new experiment.
$ws.Range("F32").Value2 = 0.3
$ws.Range("L32").Value2 = 30
$ws.Range("V32").Value2 = 0.987
$ws.Range("X32").Value2 = "log-2019-12-17T18-24-58"
$ws.Range("Y32").Value2 = "重新训练带来微小提升"
